# Update the Suzhou comic-con info workbook to match upstream gh-pages output.
#
# Summary of changes:
#  - The "苏州·明日方舟同人only （聚会）" event was cancelled: its name gets a
#    trailing "（取消）" suffix and its lowest-price column switches from a
#    numeric value to the text "不可售" ("not available for sale").
#  - A handful of "want to go" head-count numbers (column F) were refreshed
#    with newer scraped totals.
#  - These edits apply on both the "展览" (exhibitions) sheet and the
#    "全部类型" (all types) sheet, which both carry a copy of the same rows;
#    the "演出" (performances) sheet only needed one head-count refresh.

$wb = $excel.ActiveWorkbook

$cancelledNewName = "苏州·明日方舟同人only （聚会）（取消）"
$notForSale = "不可售"

# ---------------------------------------------------------------------
# Sheet "展览" (exhibitions)
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")

$ws1.Range("C3").Value = $cancelledNewName
$ws1.Range("G3").Value = $notForSale

$ws1.Range("F5").Value = 657
$ws1.Range("F6").Value = 163
$ws1.Range("F8").Value = 186
$ws1.Range("F9").Value = 368
$ws1.Range("F10").Value = 487
$ws1.Range("F11").Value = 522
$ws1.Range("F13").Value = 12195
$ws1.Range("F14").Value = 20
$ws1.Range("F15").Value = 5452

# ---------------------------------------------------------------------
# Sheet "演出" (performances)
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("演出")

$ws2.Range("F2").Value = 113

# ---------------------------------------------------------------------
# Sheet "全部类型" (all types)
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")

$ws4.Range("C3").Value = $cancelledNewName
$ws4.Range("G3").Value = $notForSale

$ws4.Range("F4").Value = 113
$ws4.Range("F7").Value = 657
$ws4.Range("F8").Value = 163
$ws4.Range("F10").Value = 186
$ws4.Range("F11").Value = 368
$ws4.Range("F12").Value = 487
$ws4.Range("F13").Value = 522
$ws4.Range("F15").Value = 12195
$ws4.Range("F17").Value = 20
$ws4.Range("F18").Value = 5452
